$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.720.34"
$ws.Range("E2").Value = "  +2.40%  "

$ws.Range("D3").Value = "'2.209.71"
$ws.Range("E3").Value = "  +1.75%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'250.94"
$ws.Range("E5").Value = "  +5.44%  "

$ws.Range("D6").Value = "'0.613"
$ws.Range("E6").Value = "  +0.90%  "

$ws.Range("D7").Value = "'74.56"
$ws.Range("E7").Value = "  +3.58%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").Value = "'0.590"
$ws.Range("E9").Value = "  +2.15%  "

$ws.Range("D10").Value = "'40.30"
$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("D11").Value = "'0.0919"
$ws.Range("E11").Value = "  +0.94%  "

$ws.Range("D12").Value = "'6.82"
$ws.Range("E12").Value = "  +1.82%  "

$ws.Range("E13").Value = "  +1.08%  "

$ws.Range("D14").Value = "'2.542.63"
$ws.Range("E14").Value = "  +1.94%  "

$ws.Range("D15").Value = "'14.45"
$ws.Range("E15").Value = "  +0.32%  "

$ws.Range("D16").Value = "'2.203.55"
$ws.Range("E16").Value = "  +2.56%  "

$ws.Range("D17").Value = "'0.782"
$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").Value = "'42.619.45"
$ws.Range("E18").Value = "  +2.81%  "

$ws.Range("E19").Value = "  +1.51%  "

$ws.Range("D20").Value = "'71.10"
$ws.Range("E20").Value = "  +1.95%  "

$ws.Range("D21").Value = "'5.93"
$ws.Range("E21").Value = "  +2.56%  "

$ws.Range("D22").Value = "'229.24"
$ws.Range("E22").Value = "  +1.29%  "

$ws.Range("D23").Value = "'2.18"
$ws.Range("E23").Value = "  +9.62%  "

$ws.Range("D24").Value = "'9.44"
$ws.Range("E24").Value = "  -5.23%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").Value = "'10.74"
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("E27").Value = "  +0.86%  "

$ws.Range("D28").Value = "'39.25"
$ws.Range("E28").Value = "  +22.45%  "

$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.21"
$ws.Range("E29").Value = "  +1.23%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  +2.47%  "

$ws.Range("D31").Value = "'170.06"
$ws.Range("E31").Value = "  -0.77%  "

$ws.Range("E32").Value = "  +1.46%  "

$ws.Range("D33").Value = "'0.0796"
$ws.Range("E33").Value = "  +3.13%  "

$ws.Range("D34").Value = "'5.21"
$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("D35").Value = "'0.121"
$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("E36").Value = "  +3.50%  "

$ws.Range("D37").Value = "'4.38"
$ws.Range("E37").Value = "  +1.90%  "

$ws.Range("D38").Value = "'0.0326"
$ws.Range("E38").Value = "  +8.95%  "

$ws.Range("D39").Value = "'12.02"
$ws.Range("E39").Value = "  -3.08%  "

$ws.Range("E40").Value = "  +1.07%  "

$ws.Range("E41").Value = "  +5.61%  "

$ws.Range("E42").Value = "  -1.32%  "

$ws.Range("D43").Value = "'58.96"
$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'102.60"
$ws.Range("E44").Value = "  +4.46%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'8.52"
$ws.Range("E45").Value = "  +2.78%  "

$ws.Range("D46").Value = "'0.475"
$ws.Range("E46").Value = "  +21.01%  "

$ws.Range("D47").Value = "'0.0979"
$ws.Range("E47").Value = "  +1.16%  "

$ws.Range("D48").Value = "'2.40"
$ws.Range("E48").Value = "  +10.39%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'1.11"
$ws.Range("E49").Value = "  +2.40%  "

$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "'1.13"
$ws.Range("E50").Value = "  +1.29%  "

$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.67"
$ws.Range("E51").Value = "  +2.14%  "
